$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 163, pushing the existing rows 163:201 down to 164:202
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new data point
$ws.Range("A163").Value2 = 9
$ws.Range("B163").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C163").Value2 = "Metropolitana"
$ws.Range("D163").Value2 = 44551
$ws.Range("E163").Value2 = 13
$ws.Range("F163").Value2 = 300000001
$ws.Range("G163").Value2 = "Rabanito"
$ws.Range("H163").Value2 = "Sin especificar"
$ws.Range("I163").Value2 = "Primera"
$ws.Range("J163").Value2 = 6100
$ws.Range("K163").Value2 = 2500
$ws.Range("L163").Value2 = 3000
$ws.Range("M163").Value2 = 2750
$ws.Range("N163").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O163").Value2 = "Provincia de Chacabuco"
$ws.Range("P163").Value2 = 28
$ws.Range("Q163").Value2 = 100
$ws.Range("R163").Value2 = "Hortaliza"
